# Updated care data (v25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21 ("Bildt ...") for the new organisation
# "Beweging 3.0 (Stichting)" which sorts alphabetically between
# "Bethanie (Stichting)" (row 20) and "Bildt (Zorgcentrum het ... Stichting)" (row 21).
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "Beweging 3.0 (Stichting)"
$ws.Range("B21").Value = "Vastgesteld"

# Remove the row for "Fundis (Stichting) (onderdeel van Welthuis)".
# After the insert above it has shifted down from row 42 to row 43.
$ws.Rows.Item(43).Delete()

# Flip status from "Voorlopig" to "Vastgesteld" for a handful of organisations.
$ws.Range("B46").Value = "Vastgesteld"
$ws.Range("B51").Value = "Vastgesteld"
$ws.Range("B63").Value = "Vastgesteld"
$ws.Range("B136").Value = "Vastgesteld"
